$d = $word.ActiveDocument

# Locate the paragraph that contains the Word field (fldChar begin/instrText/fldChar end)
# built from the TWiki "m:(...)" expression, and rewrite its runs as plain literal
# text runs using the M2Doc "{ ... }" token syntax instead of a real Word field,
# per "Updated parser to use TokenIteratorFieldRewriterSplit."

$fieldPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $s = $cand.Range.Start
    $e = $cand.Range.End
    $hit = $false
    for ($j = 1; $j -le $d.Fields.Count; $j++) {
        $fld = $d.Fields.Item($j)
        if ($fld.Code.Start -ge $s -and $fld.Code.Start -lt $e) {
            $hit = $true
        }
    }
    if ($hit) {
        $fieldPara = $cand
    }
}

if ($fieldPara -eq $null) {
    # Fallback: the second paragraph holds the TWiki field in this template.
    $fieldPara = $d.Paragraphs.Item(2)
}

$full = $fieldPara.Range
$start = $full.Start
$end = $full.End - 1
$target = $d.Range($start, $end)

$newRunsXml = '<w:r><w:t>{</w:t></w:r>' + `
  '<w:r><w:t>m</w:t></w:r>' + `
  '<w:r><w:t>:</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' + `
  "<w:r><w:t>'</w:t></w:r>" + `
  '<w:r><w:t>---+</w:t></w:r>' + `
  "<w:r><w:t xml:space=`"preserve`"> ' + </w:t></w:r>" + `
  '<w:r><w:t>self.na</w:t></w:r>' + `
  '<w:r><w:t>me).from</w:t></w:r>' + `
  '<w:r><w:t>T</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '<w:r><w:t>Wiki</w:t></w:r>' + `
  '<w:r><w:t>String()</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">}</w:t></w:r>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body><w:p>' + $newRunsXml + '</w:p></w:body>' + `
  '</w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)

Write-Host "Updated paragraph text:" $d.Paragraphs.Item(2).Range.Text
